# Insert a new data row at row 90 (shifts existing rows 90-162 down to 91-163)
# and populate it with the new reading for Terminal Hortofrutícola Agro Chillán.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("90").Insert()

$ws.Range("A90").Value2 = 7
$ws.Range("B90").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C90").Value2 = "Ñuble"
$ws.Range("D90").Value2 = 44512
$ws.Range("E90").Value2 = 16
$ws.Range("F90").Value2 = "Fruta"
$ws.Range("G90").Value2 = 100102
$ws.Range("H90").Value2 = "Cítricos"
$ws.Range("I90").Value2 = 100102004
$ws.Range("J90").Value2 = "Mandarina"
$ws.Range("K90").Value2 = "Clementina"
$ws.Range("L90").Value2 = "Primera"
$ws.Range("M90").Value2 = 120
$ws.Range("N90").Value2 = 7500
$ws.Range("O90").Value2 = 8000
$ws.Range("P90").Value2 = 7750
$ws.Range("Q90").Value2 = "`$/caja 18 kilos"
$ws.Range("R90").Value2 = "Región de O'Higgins"
$ws.Range("S90").Value2 = 431
$ws.Range("T90").Value2 = 18
